# schedule-2020.xlsx: "fixed schedule, added fol notes"
#
# The class schedule on Sheet1 had its C-column "date" shared-formula
# series re-anchored starting at row 16 (class 15, C16) so that later
# dates shift out by two days (a missed class / rescheduling), and the
# Tue/Thur weekday labels in column B for rows 24-29 were swapped to
# match the corrected dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column C: re-derive the schedule dates from row 16 (class 15) onward ---
# Row 16 used to continue the C3+7 shared-formula chain; it now restarts
# relative to the (unmoved) C15 and subsequent rows chain off each other,
# shifting every date in rows 16-29 two days later.
$ws.Range("C16").Formula = "=C15+7"
$ws.Range("C17").Formula = "=C16+2"
$ws.Range("C18").Formula = "=C16+7"
$ws.Range("C19").Formula = "=C17+7"
$ws.Range("C20").Formula = "=C18+7"
$ws.Range("C21").Formula = "=C19+7"
$ws.Range("C22").Formula = "=C20+7"
$ws.Range("C23").Formula = "=C21+7"
$ws.Range("C24").Formula = "=C22+7"
$ws.Range("C25").Formula = "=C23+7"
$ws.Range("C26").Formula = "=C24+7"
$ws.Range("C27").Formula = "=C25+7"
$ws.Range("C28").Formula = "=C26+7"
$ws.Range("C29").Formula = "=C27+7"

# --- Column B: swap the Tue/Thur weekday labels for rows 24-29 to match ---
$ws.Range("B24").Value = "Thur"
$ws.Range("B25").Value = "Tue"
$ws.Range("B26").Value = "Thur"
$ws.Range("B27").Value = "Tue"
$ws.Range("B28").Value = "Thur"
$ws.Range("B29").Value = "Tue"

# --- View state: scroll down a bit further and move the active selection ---
$ws.Activate()
$ws.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
